$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G1").Value = "Credits"
$ws.Range("G2").Value = 100
$ws.Range("G3").Value = 250
$ws.Range("G4").Value = 500

$ws.Range("G5").Select()
